# RE:L08:Slide 24 - Added a missing period.
#
# The requirements text box ("PlaceHolder 1") on slide 24 ends with the
# sentence "... previous trips" (missing a trailing full stop). This script
# appends the missing period, leaving all other runs/formatting untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)

# Find the shape named "PlaceHolder 1" which holds the requirement bullets.
$shp = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "PlaceHolder 1") {
        $shp = $candidate
        break
    }
}

$tr = $shp.TextFrame.TextRange
$full = $tr.Text

# Locate the exact run text "  previous trips" (including its leading
# whitespace, which belongs to the same run) and append a period to it in
# place, so the underlying XML run simply becomes "  previous trips.".
$needle = "  previous trips"
$idx = $full.LastIndexOf($needle)
$startPos = $idx + 1
$len = $needle.Length
$target = $tr.Characters($startPos, $len)
$target.Text = $needle + "."
